$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-9 from 45175 to 45183
for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 3).Value = 45183
}
